$wb = $excel.ActiveWorkbook

# Delete the "Childhood Giant Cell Glioblastoma" / "C114966" row from the
# "Basic Cancer Types" sheet (it's no longer available in the drop-down).
$ws = $wb.Worksheets.Item("Basic Cancer Types")
$ws.Rows.Item(4).Delete()

# Update the selection on that sheet to match the new last row.
$ws.Range("A4:XFD4").Select()

# Make "Basic Cancer Types" the active sheet/tab.
$ws.Activate()
